$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.755.29'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -5.84%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.294.04'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -6.22%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.83'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.26'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.97%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.289.27'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.06%  '
$ws.Range('E10').Value = '  -9.79%  '
$ws.Range('E11').Value = '  -6.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.33'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -8.16%  '
$ws.Range('E13').Value = '  -7.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '643.72'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('E15').Value = '  -5.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.819.28'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -6.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.12'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.723.27'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.87%  '
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.295.38'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.30%  '
$ws.Range('E21').Value = '  -8.38%  '
$ws.Range('E22').Value = '  -4.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.22'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '107.53'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +8.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.90'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -8.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.98'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -7.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.68'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -7.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.56'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.67'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -7.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -7.80%  '
$ws.Range('E31').Value = '  -8.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.28'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -7.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.08'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.104'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.800.01'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.47'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.62%  '
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '521.01'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -9.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0734'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.36'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.38%  '
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('E42').Value = '  -6.41%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '32.94'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.16%  '
$ws.Range('B44').Value = 'CoreDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.35'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -13.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.337'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -10.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0415'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.86%  '
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('E48').Value = '  -4.68%  '
$ws.Range('E49').Value = '  -9.52%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +1.29%  '
